# Auto-generated market-price refresh for Faerie_Profits workbook.
# Columns H..N = currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 797.26666  # was 840.1429000000001
$ws.Range("I33").Value = 611  # was 657
$ws.Range("K33").Value = 611  # was 657
$ws.Range("M33").Value = -382  # was -428
$ws.Range("H40").Value = 4387.778  # was 4643.0625
$ws.Range("I40").Value = 0  # was 2000
$ws.Range("J40").Value = 4387.778  # was 4819.2666
$ws.Range("K40").Value = 0  # was 2000
$ws.Range("L40").Value = 4387.778  # was 4819.2666
$ws.Range("M40").ClearContents()  # was -1825
$ws.Range("N40").Value = -4737.778  # was -5169.2666
$ws.Range("H70").Value = 9742.125  # was 10089.608
$ws.Range("J70").Value = 10263.685  # was 10736.667
$ws.Range("L70").Value = 30791.055  # was 32210.001
$ws.Range("N70").Value = -31331.055  # was -32750.001
$ws.Range("H73").Value = 9742.125  # was 10089.608
$ws.Range("J73").Value = 10263.685  # was 10736.667
$ws.Range("L73").Value = 30791.055  # was 32210.001
$ws.Range("N73").Value = -32663.055  # was -34082.001
$ws.Range("H92").Value = 1134.6562  # was 1168.7742
$ws.Range("I92").Value = 881.08  # was 914.5833
$ws.Range("K92").Value = 881.08  # was 914.5833
$ws.Range("M92").Value = 366.92  # was 333.4167
$ws.Range("H97").Value = 5221.25  # was 4252
$ws.Range("I97").Value = 2187.5  # was 1583.3334
$ws.Range("K97").Value = 6562.5  # was 4750.0002
$ws.Range("M97").Value = -6066.5  # was -4254.0002
$ws.Range("H99").Value = 486.30768  # was 484.23077
$ws.Range("I99").Value = 422.75  # was 419.375
$ws.Range("K99").Value = 1268.25  # was 1258.125
$ws.Range("M99").Value = 229.75  # was 239.875
$ws.Range("H101").Value = 1423.75  # was 1666.6666
$ws.Range("I101").Value = 1497.5  # was 1500
$ws.Range("J101").Value = 1350  # was 2000
$ws.Range("K101").Value = 4492.5  # was 4500
$ws.Range("L101").Value = 4050  # was 6000
$ws.Range("M101").Value = -2870.5  # was -2878
$ws.Range("N101").Value = -7294  # was -9244
$ws.Range("H104").Value = 1297  # was 1114
$ws.Range("I104").Value = 1256  # was 1012
$ws.Range("K104").Value = 3768  # was 3036
$ws.Range("M104").Value = -2021  # was -1289
$ws.Range("H127").Value = 868.8  # was 930.8889
$ws.Range("I127").Value = 743.1111  # was 797.25
$ws.Range("K127").Value = 2229.3333  # was 2391.75
$ws.Range("M127").Value = 2730.6667  # was 2568.25
$ws.Range("H129").Value = 45457692  # was 45457612
$ws.Range("I129").Value = 62500576  # was 58824084
$ws.Range("J129").Value = 10000  # was 11600
$ws.Range("K129").Value = 187501728  # was 176472252
$ws.Range("L129").Value = 30000  # was 34800
$ws.Range("M129").Value = -187496728  # was -176467252
$ws.Range("N129").Value = -40000  # was -44800
$ws.Range("H131").Value = 11718.417  # was 16779.715
$ws.Range("I131").Value = 13665.3  # was 19243
$ws.Range("J131").Value = 1984  # was 2000
$ws.Range("K131").Value = 40995.89999999999  # was 57729
$ws.Range("L131").Value = 5952  # was 6000
$ws.Range("M131").Value = -35955.89999999999  # was -52689
$ws.Range("N131").Value = -16032  # was -16080
$ws.Range("H132").Value = 18520668  # was 18520672
$ws.Range("I132").Value = 20410432  # was 20410436
$ws.Range("K132").Value = 61231296  # was 61231308
$ws.Range("M132").Value = -61228766  # was -61228778
$ws.Range("H138").Value = 1653.39  # was 1532.4747
$ws.Range("J138").Value = 2024.8733  # was 1859.1714
$ws.Range("L138").Value = 6074.6199  # was 5577.5142
$ws.Range("N138").Value = -16354.6199  # was -15857.5142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 4007980  # was 3341649.8
$ws.Range("J6").Value = 10000  # was 9999.75
$ws.Range("L6").Value = 10000  # was 9999.75
$ws.Range("N6").Value = -10346  # was -10345.75
$ws.Range("H61").Value = 951432.6  # was 951641.4
$ws.Range("I61").Value = 1357457.6  # was 1395152.6
$ws.Range("J61").Value = 12499.875  # was 12441.059
$ws.Range("K61").Value = 1357457.6  # was 1395152.6
$ws.Range("L61").Value = 12499.875  # was 12441.059
$ws.Range("M61").Value = -1357245.6  # was -1394940.6
$ws.Range("N61").Value = -12923.875  # was -12865.059
$ws.Range("H97").Value = 3013.7693  # was 3094.4
$ws.Range("I97").Value = 1141.6522  # was 1148.1818
$ws.Range("K97").Value = 1141.6522  # was 1148.1818
$ws.Range("M97").Value = -645.6522  # was -652.1818000000001
$ws.Range("H107").Value = 89202.09  # was 89271.73
$ws.Range("J107").Value = 89202.09  # was 89271.73
$ws.Range("L107").Value = 89202.09  # was 89271.73
$ws.Range("N107").Value = -96882.09  # was -96951.73
$ws.Range("H117").Value = 71000  # was 71333
$ws.Range("J117").Value = 71000  # was 71333
$ws.Range("L117").Value = 71000  # was 71333
$ws.Range("N117").Value = -80178  # was -80511
$ws.Range("H122").Value = 3984.0833  # was 3443.739
$ws.Range("I122").Value = 3588.0688  # was 3052.282
$ws.Range("K122").Value = 10764.2064  # was 9156.846000000001
$ws.Range("M122").Value = -8314.206399999999  # was -6706.846000000001
$ws.Range("H132").Value = 1426187.8  # was 1510039
$ws.Range("I132").Value = 1673813.1  # was 1790540.9
$ws.Range("K132").Value = 5021439.300000001  # was 5371622.699999999
$ws.Range("M132").Value = -5018909.300000001  # was -5369092.699999999
$ws.Range("H136").Value = 951432.6  # was 951641.4
$ws.Range("I136").Value = 1357457.6  # was 1395152.6
$ws.Range("J136").Value = 12499.875  # was 12441.059
$ws.Range("K136").Value = 4072372.8  # was 4185457.8
$ws.Range("L136").Value = 37499.625  # was 37323.177
$ws.Range("M136").Value = -4069822.8  # was -4182907.8
$ws.Range("N136").Value = -42599.625  # was -42423.177

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3746.52  # was 3897.6667
$ws.Range("I94").Value = 1368.7142  # was 1431.2
$ws.Range("K94").Value = 1368.7142  # was 1431.2
$ws.Range("M94").Value = -917.7141999999999  # was -980.2
$ws.Range("H99").Value = 4293.9287  # was 4312.6
$ws.Range("I99").Value = 3159.8572  # was 3170.4285
$ws.Range("J99").Value = 5428  # was 5312
$ws.Range("K99").Value = 3159.8572  # was 3170.4285
$ws.Range("L99").Value = 5428  # was 5312
$ws.Range("M99").Value = -1661.8572  # was -1672.4285
$ws.Range("N99").Value = -8424  # was -8308
$ws.Range("H105").Value = 2846.3333  # was 2901.9429
$ws.Range("I105").Value = 2086.3872  # was 2125.9333
$ws.Range("K105").Value = 2086.3872  # was 2125.9333
$ws.Range("M105").Value = -339.3872000000001  # was -378.9333000000001
$ws.Range("H134").Value = 6170.8203  # was 5615.977
$ws.Range("I134").Value = 2722.2856  # was 2446.5
$ws.Range("K134").Value = 8166.8568  # was 7339.5
$ws.Range("M134").Value = -5631.8568  # was -4804.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2059.5  # was 2060.6428
$ws.Range("I94").Value = 1457.6666  # was 1458.6666
$ws.Range("J94").Value = 2223.6365  # was 2224.818
$ws.Range("K94").Value = 1457.6666  # was 1458.6666
$ws.Range("L94").Value = 2223.6365  # was 2224.818
$ws.Range("M94").Value = -1006.6666  # was -1007.6666
$ws.Range("N94").Value = -3125.6365  # was -3126.818

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 141.58621  # was 141.7931
$ws.Range("I2").Value = 146.1  # was 155.44444
$ws.Range("J2").Value = 131.55556  # was 119.454544
$ws.Range("K2").Value = 876.5999999999999  # was 932.6666399999999
$ws.Range("L2").Value = 789.3333600000001  # was 716.727264
$ws.Range("M2").Value = -763.5999999999999  # was -819.6666399999999
$ws.Range("N2").Value = -1015.33336  # was -942.727264
$ws.Range("H113").Value = 1625  # was 1566.3462
$ws.Range("I113").Value = 880  # was 490
$ws.Range("K113").Value = 2640  # was 1470
$ws.Range("M113").Value = -470  # was 700
$ws.Range("H134").Value = 7045  # was 7105.8
$ws.Range("J134").Value = 7979.3335  # was 10000
$ws.Range("L134").Value = 23938.0005  # was 30000
$ws.Range("N134").Value = -34078.00049999999  # was -40140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8500  # was 4188.1665
$ws.Range("I122").Value = 0  # was 3000
$ws.Range("J122").Value = 8500  # was 4425.8
$ws.Range("K122").Value = 0  # was 9000
$ws.Range("L122").Value = 25500  # was 13277.4
$ws.Range("M122").ClearContents()  # was -6550
$ws.Range("N122").Value = -30400  # was -18177.4
$ws.Range("H123").Value = 43149.25  # was 44133
$ws.Range("J123").Value = 43149.25  # was 44133
$ws.Range("L123").Value = 43149.25  # was 44133
$ws.Range("N123").Value = -48049.25  # was -49033
$ws.Range("H132").Value = 5558705  # was 5750435.5
$ws.Range("I132").Value = 7410425.5  # was 7939680
$ws.Range("J132").Value = 3543.2  # was 3668.625
$ws.Range("K132").Value = 22231276.5  # was 23819040
$ws.Range("L132").Value = 10629.6  # was 11005.875
$ws.Range("M132").Value = -22228746.5  # was -23816510
$ws.Range("N132").Value = -15689.6  # was -16065.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 10680.421  # was 11225.223
$ws.Range("I46").Value = 7258  # was 10450
$ws.Range("K46").Value = 7258  # was 10450
$ws.Range("M46").Value = -7070  # was -10262
$ws.Range("H63").Value = 83637  # was 85701.664
$ws.Range("J63").Value = 83637  # was 85701.664
$ws.Range("L63").Value = 83637  # was 85701.664
$ws.Range("N63").Value = -85135  # was -87199.664
$ws.Range("H64").Value = 0  # was 20000
$ws.Range("J64").Value = 0  # was 20000
$ws.Range("L64").Value = 0  # was 20000
$ws.Range("N64").ClearContents()  # was -20450
$ws.Range("H66").Value = 83637  # was 85701.664
$ws.Range("J66").Value = 83637  # was 85701.664
$ws.Range("L66").Value = 250911  # was 257104.992
$ws.Range("N66").Value = -258399  # was -264592.992
$ws.Range("H67").Value = 0  # was 20000
$ws.Range("J67").Value = 0  # was 20000
$ws.Range("L67").Value = 0  # was 20000
$ws.Range("N67").ClearContents()  # was -21560
$ws.Range("H101").Value = 50000  # was 70000
$ws.Range("J101").Value = 50000  # was 70000
$ws.Range("L101").Value = 50000  # was 70000
$ws.Range("N101").Value = -56490  # was -76490
$ws.Range("H115").Value = 90000  # was 0
$ws.Range("J115").Value = 90000  # was 0
$ws.Range("L115").Value = 90000  # was 0
$ws.Range("N115").Value = -92350  # was (absent)
$ws.Range("H118").Value = 90000  # was 0
$ws.Range("J118").Value = 90000  # was 0
$ws.Range("L118").Value = 90000  # was 0
$ws.Range("N118").Value = -93314  # was (absent)
$ws.Range("H122").Value = 7743.8184  # was 7124.4
$ws.Range("I122").Value = 7132.1665  # was 6588.5557
$ws.Range("J122").Value = 8477.799999999999  # was 7928.1665
$ws.Range("K122").Value = 21396.4995  # was 19765.6671
$ws.Range("L122").Value = 25433.4  # was 23784.4995
$ws.Range("M122").Value = -18946.4995  # was -17315.6671
$ws.Range("N122").Value = -30333.4  # was -28684.4995
$ws.Range("H132").Value = 3131.1667  # was 3122.0945
$ws.Range("J132").Value = 3148.875  # was 3121.6924
$ws.Range("L132").Value = 9446.625  # was 9365.0772
$ws.Range("N132").Value = -14506.625  # was -14425.0772
$ws.Range("H136").Value = 5150.724  # was 5192.7144
$ws.Range("I136").Value = 4685.478  # was 4717.773
$ws.Range("K136").Value = 14056.434  # was 14153.319
$ws.Range("M136").Value = -11506.434  # was -11603.319

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 101204  # was 99738.164
$ws.Range("J46").Value = 101204  # was 99738.164
$ws.Range("L46").Value = 101204  # was 99738.164
$ws.Range("N46").Value = -101666  # was -100200.164
$ws.Range("H100").Value = 821.6  # was 829.5806
$ws.Range("I100").Value = 716.6  # was 730.1539
$ws.Range("K100").Value = 1433.2  # was 1460.3078
$ws.Range("M100").Value = -892.2  # was -919.3078
$ws.Range("H116").Value = 58340  # was 0
$ws.Range("J116").Value = 58340  # was 0
$ws.Range("L116").Value = 58340  # was 0
$ws.Range("N116").Value = -67518  # was (absent)
$ws.Range("H122").Value = 2552.3157  # was 2246.348
$ws.Range("I122").Value = 2383.0833  # was 2156
$ws.Range("J122").Value = 2842.4285  # was 2386.889
$ws.Range("K122").Value = 7149.249899999999  # was 6468
$ws.Range("L122").Value = 8527.2855  # was 7160.667
$ws.Range("M122").Value = -4699.249899999999  # was -4018
$ws.Range("N122").Value = -13427.2855  # was -12060.667
$ws.Range("H132").Value = 2199  # was 2261.5
$ws.Range("I132").Value = 2073.5  # was 2198.3333
$ws.Range("K132").Value = 6220.5  # was 6594.999899999999
$ws.Range("M132").Value = -3690.5  # was -4064.999899999999
$ws.Range("H134").Value = 101204  # was 99738.164
$ws.Range("J134").Value = 101204  # was 99738.164
$ws.Range("L134").Value = 303612  # was 299214.492
$ws.Range("N134").Value = -308682  # was -304284.492
